# "desenat bradut in consola"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix typo in student name (B7 shared string)
$ws.Range("B7").Value = "Cătălina Mădălina Paca"

# 2) Add new attendance marks (K column, week 9) for several students
$ws.Range("K7").Value = 2
$ws.Range("K10").Value = 2
$ws.Range("K13").Value = 2
$ws.Range("K15").Value = 2
$ws.Range("K20").Value = 2
$ws.Range("K22").Value = 2

# 3) Re-enter the attendance-total formulas so Q3:Q22 become one shared
#    formula group (Q23:Q52 was already a shared group and keeps being one)
$ws.Range("Q3:Q22").Formula = "=SUM(C3:P3)"

# 4) Merge the two conditional-formatting rules on column Q into a single
#    rule spanning Q3:Q52, keeping the formatting (dxf) that was used by
#    the Q22-only rule.
$bigRule = $ws.Range("Q3:Q21,Q23:Q52").FormatConditions.Item(1)
$bigRule.Delete()
$keepRule = $ws.Range("Q22").FormatConditions.Item(1)
$keepRule.ModifyAppliesToRange($ws.Range("Q3:Q52"))

# 5) Update the active cell selection on the sheet
$ws.Range("B7").Select()
